$d = $word.ActiveDocument

# 1. Move the "_GoBack" bookmark from the "Queue/Stack Practice" paragraph
#    to the empty paragraph that follows the valgrind instructions
#    (just before "int main();").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Re-type a handful of paragraphs whose text spans spell-check
#    (proofErr) boundaries so Word collapses the split runs back into a
#    single run per paragraph (clearing the now-stale proofErr markers),
#    and fix the valgrind flag dashes from an en-dash to a literal "--".
$d.Content.Find.Execute("bool isPalindrome(char *);", $false, $false, $false, $false, $false, $true, 1, $false, "bool isPalindrome(char *);", 2) | Out-Null

$d.Content.Find.Execute("Returns a bool.", $false, $false, $false, $false, $false, $true, 1, $false, "Returns a bool.", 2) | Out-Null

$d.Content.Find.Execute("The dequeue and pop functions should free the node memory as we saw in class today.  To check for memory leaks, use valgrind by running", $false, $false, $false, $false, $false, $true, 1, $false, "The dequeue and pop functions should free the node memory as we saw in class today.  To check for memory leaks, use valgrind by running", 2) | Out-Null

$d.Content.Find.Execute("%valgrind –leak-check=full <program name>", $false, $false, $false, $false, $false, $true, 1, $false, "%valgrind --leak-check=full <program name>", 2) | Out-Null

$d.Content.Find.Execute("int main();", $false, $false, $false, $false, $false, $true, 1, $false, "int main();", 2) | Out-Null

$d.Content.Find.Execute("Call the isPalindrome() function with the word as the parameter and compare the output of the function to the 0/1 flag from the file.", $false, $false, $false, $false, $false, $true, 1, $false, "Call the isPalindrome() function with the word as the parameter and compare the output of the function to the 0/1 flag from the file.", 2) | Out-Null

# 3. Add the bookmark back on the empty paragraph right after the
#    "%valgrind --leak-check=full <program name>" line.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("leak-check")) {
        $target = $p.Next()
        $rng = $target.Range
        $rng.End = $rng.End - 1
        $d.Bookmarks.Add("_GoBack", $rng)
        break
    }
}
